$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the numbering sequence in column A (row 7 was blank, now "4")
$ws.Range("A7").Value = 4

# Move the active selection, matching the saved cursor position
$ws.Range("I12").Select()
